# Regenerate save_data: recompute column G ("K", formerly derived from Strike#)
# for each data row (rows 2-67) using the newly regenerated std/mean based
# s_vals calculation, and write the resulting K values back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2 through 67, in row order.
$newKValues = @(
    3, 0, 1, 3, 3, 0, 3, 0, 0, 1, 2, 1, 3, 1, 1, 3, 3, 0, 1, 1,
    1, 1, 0, 2, 0, 2, 0, 0, 1, 1, 1, 1, 0, 1, 2, 2, 1, 1, 3, 1,
    3, 1, 3, 2, 2, 4, 0, 1, 2, 0, 1, 2, 1, 1, 3, 1, 0, 2, 1, 1,
    2, 1, 0, 1, 2, 2
)

$startRow = 2
for ($i = 0; $i -lt $newKValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newKValues[$i]
}
